# Auto-generated edit script: apply numeric corrections to the
# per-leve profit/price columns (H:N) across all 8 crafting-job sheets,
# per commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 181
$ws.Range("I4").Value = 149.2
$ws.Range("J4").Value = 340
$ws.Range("K4").Value = 149.2
$ws.Range("L4").Value = 340
$ws.Range("M4").Value = -35.19999999999999
$ws.Range("N4").Value = -568
$ws.Range("H127").Value = 2208.5789
$ws.Range("J127").Value = 2528.9375
$ws.Range("L127").Value = 7586.8125
$ws.Range("N127").Value = -17506.8125
$ws.Range("H129").Value = 786.75
$ws.Range("J129").Value = 888.2593000000001
$ws.Range("L129").Value = 2664.7779
$ws.Range("N129").Value = -12664.7779
$ws.Range("H132").Value = 3379.3572
$ws.Range("I132").Value = 3379.3572
$ws.Range("K132").Value = 10138.0716
$ws.Range("M132").Value = -7608.071599999999
$ws.Range("H137").Value = 710.46155
$ws.Range("I137").Value = 656.34784
$ws.Range("J137").Value = 1125.3334
$ws.Range("K137").Value = 1969.04352
$ws.Range("L137").Value = 3376.0002
$ws.Range("M137").Value = 580.9564799999998
$ws.Range("N137").Value = -8476.0002

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 779.7
$ws.Range("I2").Value = 817.725
$ws.Range("J2").Value = 627.6
$ws.Range("K2").Value = 817.725
$ws.Range("L2").Value = 627.6
$ws.Range("M2").Value = -704.725
$ws.Range("N2").Value = -853.6
$ws.Range("H32").Value = 6427.6943
$ws.Range("I32").Value = 5208.8804
$ws.Range("J32").Value = 22759.8
$ws.Range("K32").Value = 5208.8804
$ws.Range("L32").Value = 22759.8
$ws.Range("M32").Value = -4921.8804
$ws.Range("N32").Value = -23333.8
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H110").Value = 39112.23
$ws.Range("I110").Value = 46151
$ws.Range("J110").Value = 399
$ws.Range("K110").Value = 46151
$ws.Range("L110").Value = 399
$ws.Range("M110").Value = -44106
$ws.Range("N110").Value = -4489
$ws.Range("H116").Value = 779.7
$ws.Range("I116").Value = 817.725
$ws.Range("J116").Value = 627.6
$ws.Range("K116").Value = 817.725
$ws.Range("L116").Value = 627.6
$ws.Range("M116").Value = 1476.275
$ws.Range("N116").Value = -5215.6
$ws.Range("H132").Value = 1830.3715
$ws.Range("I132").Value = 1195.8636
$ws.Range("J132").Value = 2904.1538
$ws.Range("K132").Value = 3587.5908
$ws.Range("L132").Value = 8712.4614
$ws.Range("M132").Value = -1057.5908
$ws.Range("N132").Value = -13772.4614
$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -70140
$ws.Range("H139").Value = 38646.848
$ws.Range("J139").Value = 38646.848
$ws.Range("L139").Value = 38646.848
$ws.Range("N139").Value = -48926.848

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 779.7
$ws.Range("I3").Value = 817.725
$ws.Range("J3").Value = 627.6
$ws.Range("K3").Value = 817.725
$ws.Range("L3").Value = 627.6
$ws.Range("M3").Value = -703.725
$ws.Range("N3").Value = -855.6
$ws.Range("H86").Value = 8650
$ws.Range("I86").Value = 8866.666999999999
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 8866.666999999999
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = -7743.666999999999
$ws.Range("N86").Value = -10246
$ws.Range("H89").Value = 8650
$ws.Range("I89").Value = 8866.666999999999
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 44333.335
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = -38717.335
$ws.Range("N89").Value = -51232
$ws.Range("H99").Value = 1843.75
$ws.Range("I99").Value = 1700
$ws.Range("J99").Value = 1930
$ws.Range("K99").Value = 1700
$ws.Range("L99").Value = 1930
$ws.Range("M99").Value = -202
$ws.Range("N99").Value = -4926

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3703.6
$ws.Range("I122").Value = 2384.8
$ws.Range("J122").Value = 5022.4
$ws.Range("K122").Value = 7154.400000000001
$ws.Range("L122").Value = 15067.2
$ws.Range("M122").Value = -4704.400000000001
$ws.Range("N122").Value = -19967.2

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 908709.2
$ws.Range("J107").Value = 1047664.25
$ws.Range("L107").Value = 3142992.75
$ws.Range("N107").Value = -3146832.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1100.4736
$ws.Range("I2").Value = 1157.3334
$ws.Range("J2").Value = 1049.3
$ws.Range("K2").Value = 1157.3334
$ws.Range("L2").Value = 1049.3
$ws.Range("M2").Value = -1044.3334
$ws.Range("N2").Value = -1275.3
$ws.Range("H70").Value = 5649.8
$ws.Range("I70").Value = 4732.722
$ws.Range("J70").Value = 6620.8237
$ws.Range("K70").Value = 4732.722
$ws.Range("L70").Value = 6620.8237
$ws.Range("M70").Value = -4462.722
$ws.Range("N70").Value = -7160.8237
$ws.Range("H73").Value = 5649.8
$ws.Range("I73").Value = 4732.722
$ws.Range("J73").Value = 6620.8237
$ws.Range("K73").Value = 4732.722
$ws.Range("L73").Value = 6620.8237
$ws.Range("M73").Value = -3796.722
$ws.Range("N73").Value = -8492.823700000001
$ws.Range("H102").Value = 2587.652
$ws.Range("I102").Value = 2723.2307
$ws.Range("J102").Value = 2411.4
$ws.Range("K102").Value = 2723.2307
$ws.Range("L102").Value = 2411.4
$ws.Range("M102").Value = -1101.2307
$ws.Range("N102").Value = -5655.4
$ws.Range("H122").Value = 2570
$ws.Range("I122").Value = 2008.8
$ws.Range("J122").Value = 3817.111
$ws.Range("K122").Value = 6026.4
$ws.Range("L122").Value = 11451.333
$ws.Range("M122").Value = -3576.4
$ws.Range("N122").Value = -16351.333

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2432.6924
$ws.Range("I68").Value = 1875.5
$ws.Range("J68").Value = 2680.3333
$ws.Range("K68").Value = 1875.5
$ws.Range("L68").Value = 2680.3333
$ws.Range("M68").Value = -1126.5
$ws.Range("N68").Value = -4178.3333
$ws.Range("H71").Value = 2432.6924
$ws.Range("I71").Value = 1875.5
$ws.Range("J71").Value = 2680.3333
$ws.Range("K71").Value = 9377.5
$ws.Range("L71").Value = 13401.6665
$ws.Range("M71").Value = -5633.5
$ws.Range("N71").Value = -20889.6665
$ws.Range("H100").Value = 1763.8667
$ws.Range("I100").Value = 1570.7273
$ws.Range("J100").Value = 2295
$ws.Range("K100").Value = 1570.7273
$ws.Range("L100").Value = 2295
$ws.Range("M100").Value = -1029.7273
$ws.Range("N100").Value = -3377
$ws.Range("H132").Value = 7777.478
$ws.Range("I132").Value = 6193.5273
$ws.Range("J132").Value = 14000.143
$ws.Range("K132").Value = 18580.5819
$ws.Range("L132").Value = 42000.429
$ws.Range("M132").Value = -16050.5819
$ws.Range("N132").Value = -47060.429
$ws.Range("H136").Value = 13933377
$ws.Range("I136").Value = 65419.938
$ws.Range("J136").Value = 41669292
$ws.Range("K136").Value = 196259.814
$ws.Range("L136").Value = 125007876
$ws.Range("M136").Value = -193709.814
$ws.Range("N136").Value = -125012976
$ws.Range("H139").Value = 44245.453
$ws.Range("J139").Value = 44245.453
$ws.Range("L139").Value = 44245.453
$ws.Range("N139").Value = -54525.453

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 113967.22
$ws.Range("I122").Value = 169066.67
$ws.Range("J122").Value = 3768.3333
$ws.Range("K122").Value = 507200.01
$ws.Range("L122").Value = 11304.9999
$ws.Range("M122").Value = -504750.01
$ws.Range("N122").Value = -16204.9999
$ws.Range("H136").Value = 9541537
$ws.Range("I136").Value = 17167352
$ws.Range("J136").Value = 9268.125
$ws.Range("K136").Value = 51502056
$ws.Range("L136").Value = 27804.375
$ws.Range("M136").Value = -51499506
$ws.Range("N136").Value = -32904.375

